# feat: add 2022-Q1 data
#
# The workbook has 6 sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计.
# We need to:
#   1. Turn the existing "总计" sheet (slot 6) into the new "2022-Q1" sheet
#      (per-fund holdings for the new quarter) - this keeps its internal slot.
#   2. Add a brand-new "总计" sheet at the end with the historical totals
#      table, including a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet (index 6) into "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)

# Before we touch any values, stamp the header-style (s=2) + index-column
# style (s=2) onto the additional cells we are about to need, by copying
# from cells that already carry that style within this same sheet.
$q1.Range("B1").Copy($q1.Range("E1:H1"))
$q1.Range("A2").Copy($q1.Range("A7:A10"))

$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G must stay text (fund codes have leading zeros, and
# the decimal figures are stored as text in the source data) - force text
# number format before writing so the engine doesn't coerce to numeric.
$q1.Range("B2:B10").NumberFormat = "@"
$q1.Range("D2:G10").NumberFormat = "@"

$fundRows = @(
    @("003751", "万家瑞隆混合", "27.84", "86.40", "2.51", "0.6988", 9),
    @("005106", "银华农业产业股票", "13.24", "93.41", "4.57", "0.6051", 9),
    @("001940", "农银汇理现代农业加灵活配置混合", "1.33", "64.68", "7.38", "0.0982", 1),
    @("009169", "湘财长兴灵活配置混合A", "1.16", "85.40", "2.85", "0.0331", 9),
    @("009500", "国寿安保高股息混合A", "0.96", "73.47", "3.06", "0.0294", 9),
    @("164401", "前海开源中证健康产业指数", "2.13", "94.15", "1.19", "0.0253", 9),
    @("008116", "银华沪深股通精选混合", "0.55", "88.86", "2.58", "0.0142", 9),
    @("009170", "湘财长兴灵活配置混合C", "0.46", "85.40", "2.85", "0.0131", 9),
    @("009501", "国寿安保高股息混合C", "0.03", "73.47", "3.06", "0.0009", 9)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: add the new "总计" sheet at the end with the totals table
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Borrow the s=2 style (header + index column) from the 2022-Q1 sheet so the
# new sheet matches the house style used across every other sheet.
$q1.Range("B1").Copy($total.Range("B1:D1"))
$q1.Range("A2").Copy($total.Range("A2:A7"))

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 9, 1.52),
    @("2021-Q4", 3, 0.67),
    @("2021-Q3", 5, 0.7),
    @("2021-Q2", 3, 1.19),
    @("2021-Q1", 5, 2.7),
    @("2020-Q4", 12, 4.5)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

Write-Host "done"
